{"js": "// Insert a new leading paragraph \"Papers and Datasets: <link>\" followed\n// by an empty paragraph, right before the existing first paragraph\n// (which starts with the \"fast.ai\" hyperlink).\n\nconst body = context.document.body;\nbody.paragraphs.load(\"items\");\nawait context.sync();\n\nconst firstPara = body.paragraphs.items[0];\n\n// Create the new paragraph before the current first paragraph.\nconst newPara = firstPara.insertParagraph(\"Papers and Datasets: \", Word.InsertLocation.before);\n\n// Append the hyperlink text and turn it into a real hyperlink whose\n// display text is the URL itself (matches the target markup).\nconst url = \"https://github.com/sebastianruder/NLP-progress\";\nconst linkRange = newPara.insertText(url, Word.InsertLocation.end);\nlinkRange.hyperlink = url;\n\n// Insert a blank paragraph between the new paragraph and the original\n// first paragraph.\nfirstPara.insertParagraph(\"\", Word.InsertLocation.before);\n\nawait context.sync();\n", "ps1": "# Insert a new leading paragraph \"Papers and Datasets: <link>\" followed\n# by an empty paragraph, right before the existing first paragraph\n# (which starts with the \"fast.ai\" hyperlink).\n\n$d = $word.ActiveDocument\n\n# Create a brand new, empty paragraph right before the current first\n# paragraph.\n$firstPara = $d.Paragraphs.Item(1)\n$firstRange = $firstPara.Range\n$firstRange.SetRange(0, 0)\n$firstRange.InsertParagraphBefore()\n\n# That new paragraph is now Paragraph 1 (still empty). Turn its whole\n# (empty) range into a hyperlink whose visible text is the URL itself.\n$newPara = $d.Paragraphs.Item(1)\n$newRange = $newPara.Range\n$url = \"https://github.com/sebastianruder/NLP-progress\"\n$d.Hyperlinks.Add($newRange, $url, [Type]::Missing, [Type]::Missing, $url)\n\n# Prepend the plain-text label before the hyperlink we just created.\n$newPara2 = $d.Paragraphs.Item(1)\n$newPara2.Range.InsertBefore(\"Papers and Datasets: \")\n\n# Insert a blank paragraph right after this one, separating it from the\n# original first paragraph.\n$para1 = $d.Paragraphs.Item(1)\n$para1.Range.InsertParagraphAfter()\n"}
